$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 and 3 currently hold the Token and WorkspaceID concatenated into a
# single string in column B (with column C left empty). Split them apart so
# the data rows match the two-column Token / WorkspaceID header layout
# already present in row 1 (B1="Token", C1="WorkspaceID").
$ws.Range("B2").Value = "ODhVGFcsg4tFPfntsKygHF3thH9WPmUfAUL2d2rh"
$ws.Range("C2").Value = "cdc95a97-6e85-4c31-9d28-230018d40671"

$ws.Range("B3").Value = "ODhVGFcsg4tFPfntsKygHF3thH9WPmUfAUL2d2rh"
$ws.Range("C3").Value = "cdc95a97-6e85-4c31-9d28-230018d40671"
